# "vul data aan met werkzaamheidsgraad"
# Add a new column L "werkzaamheidsgraad_2016" with one value per region row
# on Sheet2, formatted with a custom "0.0" number format, and select P10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Header for the new column
$ws.Range("L1").Value = "werkzaamheidsgraad_2016"

# New "werkzaamheidsgraad_2016" values, one per region row (2-12)
$values = @{
    2  = 69.2
    3  = 59.8
    4  = 60
    5  = 69.7
    6  = 61.1
    7  = 68.2
    8  = 64.2
    9  = 74
    10 = 73.9
    11 = 69.3
    12 = 74
}

foreach ($row in 2..12) {
    $cell = $ws.Range("L$row")
    $cell.Value = $values[$row]
    $cell.NumberFormat = "0.0"
}

# Restore the active selection like in the edited workbook
[void]$ws.Range("P10").Select()

# Match the page orientation recorded alongside the new column
$ws.PageSetup.Orientation = 1
